$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 (Time 43494, MSFT): Type SELL -> BUY, Amount 2.1852 -> 33.5321
$ws.Range("C5").Value = "BUY"
$ws.Range("D5").Value = 33.5321

# Row 18 (Time 43507, MSFT): Amount 1.8239000000000001 -> 8.2939000000000007
$ws.Range("D18").Value = 8.2939000000000007

# Row 28 (Time 43517, AAPL): Amount 0.7 -> 3.7
$ws.Range("D28").Value = 3.7

# Update selection to reflect the final active cell used when the file was saved
$ws.Range("H16").Select()
